$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF"), matching the bold/bordered
# --- header style already used by A1:H1 (copy format from H1, which
# --- carries cellXf style index 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-32: new numeric values for columns I (I0) and J (IF).
$data = @(
  @(2, 7, 7),
  @(3, 8, 8),
  @(4, 8, 8),
  @(5, 7, 7),
  @(6, 7, 7),
  @(7, 7, 7),
  @(8, 8, 8),
  @(9, 6, 7),
  @(10, 6, 6),
  @(11, 6, 6),
  @(12, 6, 6),
  @(13, 7, 7),
  @(14, 7, 7),
  @(15, 9, 9),
  @(16, 6, 6),
  @(17, 7, 7),
  @(18, 7, 7),
  @(19, 3, 4),
  @(20, 7, 8),
  @(21, 7, 7),
  @(22, 8, 8),
  @(23, 7, 7),
  @(24, 9, 9),
  @(25, 6, 6),
  @(26, 5, 5),
  @(27, 9, 9),
  @(28, 7, 7),
  @(29, 5, 5),
  @(30, 3, 3),
  @(31, 5, 5),
  @(32, 1, 1)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 9).Value = $row[1]
  $ws.Cells.Item($r, 10).Value = $row[2]
}
